# Add a new "canonical SMILES" column (D) that duplicates the values from the
# existing "canonical isomeric SMILES" column (C), as part of adding updated
# microstates from SAMPL6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D2").Value = "canonical SMILES"

# Copy each row's canonical isomeric SMILES (column C) into the new
# canonical SMILES column (column D)
for ($row = 3; $row -le 10; $row++) {
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value()
}

# Match the new column's width from the target workbook
$ws.Columns.Item(4).ColumnWidth = 36
